# FIX #12091 - update template with [attachmentRecipient. instead of [recipient.
# Also refresh the generated date field value and a couple of cosmetic
# table-cell-margin / decorative-line tweaks that came along with the
# template re-save.

$d = $word.ActiveDocument

# --- 1) Table cell left margins: 128 -> 133 dxa on both tables --------------
$t1 = $d.Tables.Item(1)
$t1.LeftPadding = 6.65   # 133 dxa = 133/20 pt
$t2 = $d.Tables.Item(2)
$t2.LeftPadding = 6.65   # 133 dxa = 133/20 pt

# --- 2) [recipient.xxx] -> [attachmentRecipient.xxx] placeholders ----------
$d.Content.Find.Execute(
    "[recipient.postal_address;strconv=no]", $false, $false, $false,
    $false, $false, $true, 1, $false,
    "[attachmentRecipient.postal_address;strconv=no]", 2)

$d.Content.Find.Execute(
    "[recipient.civility] [recipient.lastname],", $false, $false, $false,
    $false, $false, $true, 1, $false,
    "[attachmentRecipient.civility] [attachmentRecipient.lastname],", 2)

$d.Content.Find.Execute(
    "Veuillez agréer, [recipient.civility], l’expression de nos salutations distinguées.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Veuillez agréer, [attachmentRecipient.civility], l’expression de nos salutations distinguées.", 2)

# --- 3) Refresh the cached TIME field text ----------------------------------
$d.Content.Find.Execute(
    "09/12/2019", $false, $false, $false, $false, $false, $true, 1, $false,
    "02/01/2020", 2)

# --- 4) Nudge the decorative header line's stored size ----------------------
$sec = $d.Sections.Item(1)
$hdr = $sec.Headers.Item(1)
$lineShape = $hdr.Shapes.Item(2)
$lineShape.Width = 543.89999
$lineShape.Height = 0.34999
